$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data0 = New-Object 'object[,]' 24,2
$data0[0,0] = 0.7745651815740189
$data0[0,1] = 0.2107446106586224
$data0[1,0] = 0.6760686717638578
$data0[1,1] = 0.1886257031920877
$data0[2,0] = 0.6153414176308161
$data0[2,1] = 0.1749848748368379
$data0[3,0] = 0.5905335357662693
$data0[3,1] = 0.169411540465461
$data0[4,0] = 0.5864105738365311
$data0[4,1] = 0.1684852244008823
$data0[5,0] = 0.615007094540772
$data0[5,1] = 0.1749097693794965
$data0[6,0] = 0.7406565370809517
$data0[6,1] = 0.2031306606435805
$data0[7,0] = 0.9850026100980358
$data0[7,1] = 0.2579814587155624
$data0[8,0] = 1.163198263354502
$data0[8,1] = 0.2979634218427805
$data0[9,0] = 1.243961797043823
$data0[9,1] = 0.3160799059828889
$data0[10,0] = 1.274500336328515
$data0[10,1] = 0.3229294875474693
$data0[11,0] = 1.267925341044645
$data0[11,1] = 0.3214547923653583
$data0[12,0] = 1.246475131043269
$data0[12,1] = 0.3166436429728492
$data0[13,0] = 1.233330359346553
$data0[13,1] = 0.3136952632483485
$data0[14,0] = 1.157914003698693
$data0[14,1] = 0.2967779890596205
$data0[15,0] = 1.111570668137119
$data0[15,1] = 0.2863811488266208
$data0[16,0] = 1.0848872191429
$data0[16,1] = 0.2803944572364401
$data0[17,0] = 1.075847913744269
$data0[17,1] = 0.278366328811785
$data0[18,0] = 1.116506905983897
$data0[18,1] = 0.2874886073704488
$data0[19,0] = 1.25277681037943
$data0[19,1] = 0.3180570891959462
$data0[20,0] = 1.341574721066991
$data0[20,1] = 0.3379725629026211
$data0[21,0] = 1.294206241795621
$data0[21,1] = 0.3273491948039862
$data0[22,0] = 1.114275357453494
$data0[22,1] = 0.2869879546394429
$data0[23,0] = 0.9191279696522088
$data0[23,1] = 0.2431972319100169
$ws.Range("B2:C25").Value = $data0

$data1 = New-Object 'object[,]' 24,6
$data1[0,0] = 0.6278702139078263
$data1[0,1] = 1.725261838504053
$data1[0,2] = 0.2032907514144355
$data1[0,3] = 0.3882364568154699
$data1[0,4] = 0.2692454659551169
$data1[0,5] = 0.02373654028034622
$data1[1,0] = 0.6234773292514859
$data1[1,1] = 1.721399550507343
$data1[1,2] = 0.2059384260378536
$data1[1,3] = 0.3933223059687663
$data1[1,4] = 0.2761921925181809
$data1[1,5] = 0.02356252043098728
$data1[2,0] = 0.6210767214966921
$data1[2,1] = 1.720174844028577
$data1[2,2] = 0.207834093455908
$data1[2,3] = 0.3966954844843116
$data1[2,4] = 0.2807237040861592
$data1[2,5] = 0.02346550770905509
$data1[3,0] = 0.620173081797077
$data1[3,1] = 1.719963860827221
$data1[3,2] = 0.2086742126080097
$data1[3,3] = 0.3981330096883653
$data1[3,4] = 0.2826371095690945
$data1[3,5] = 0.02342845941618421
$data1[4,0] = 0.6200275414740233
$data1[4,1] = 1.719946219670817
$data1[4,2] = 0.2088177898403245
$data1[4,3] = 0.3983755088775638
$data1[4,4] = 0.282958857642539
$data1[4,5] = 0.02342245802286236
$data1[5,0] = 0.6210642324863827
$data1[5,1] = 1.720170832506099
$data1[5,2] = 0.2078451501679552
$data1[5,3] = 0.3967146167266051
$data1[5,4] = 0.2807492387719561
$data1[5,5] = 0.02346499798468038
$data1[6,0] = 0.6262939970799906
$data1[6,1] = 1.7236919463856
$data1[6,2] = 0.204147471910062
$data1[6,3] = 0.3899380481547396
$data1[6,4] = 0.2715853543800755
$data1[6,5] = 0.02367450367178492
$data1[7,0] = 0.6389030524691393
$data1[7,1] = 1.73971087403649
$data1[7,2] = 0.1990507812731295
$data1[7,3] = 0.3786385543754278
$data1[7,4] = 0.2557339387256397
$data1[7,5] = 0.0241629198062796
$data1[8,0] = 0.6496033412404714
$data1[8,1] = 1.7570613962274
$data1[8,2] = 0.1966366412213958
$data1[8,3] = 0.3715524586115038
$data1[8,4] = 0.2453885771586828
$data1[8,5] = 0.02456849209879053
$data1[9,0] = 0.654783577684384
$data1[9,1] = 1.766172256843021
$data1[9,2] = 0.1958309025314051
$data1[9,3] = 0.3685933958786407
$data1[9,4] = 0.2409665214714698
$data1[9,5] = 0.02476303276501568
$data1[10,0] = 0.6567901559399161
$data1[10,1] = 1.769797827979815
$data1[10,2] = 0.1955681368462692
$data1[10,3] = 0.3675109574959663
$data1[10,4] = 0.2393330245773306
$data1[10,5] = 0.0248381340449555
$data1[11,0] = 0.6563560052101778
$data1[11,1] = 1.769009186244489
$data1[11,2] = 0.1956228399470916
$data1[11,3] = 0.3677423846181043
$data1[11,4] = 0.2396829991511336
$data1[11,5] = 0.02482189608868879
$data1[12,0] = 0.6549477596270421
$data1[12,1] = 1.766467015842849
$data1[12,2] = 0.1958084342052118
$data1[12,3] = 0.3685035791598423
$data1[12,4] = 0.2408313091383985
$data1[12,5] = 0.02476918273835693
$data1[13,0] = 0.6540910198609424
$data1[13,1] = 1.764932726468132
$data1[13,2] = 0.1959276399671666
$data1[13,3] = 0.3689747958504626
$data1[13,4] = 0.2415400322529511
$data1[13,5] = 0.02473708056728796
$data1[14,0] = 0.6492710896714868
$data1[14,1] = 1.756490518511328
$data1[14,2] = 0.1966952084579674
$data1[14,3] = 0.3717511673558036
$data1[14,4] = 0.2456833016669968
$data1[14,5] = 0.02455597954385524
$data1[15,0] = 0.6463942744263349
$data1[15,1] = 1.751623701304496
$data1[15,2] = 0.1972412039981961
$data1[15,3] = 0.373522159115474
$data1[15,4] = 0.2482979610005094
$data1[15,5] = 0.02444744428646572
$data1[16,0] = 0.644769032912123
$data1[16,1] = 1.748939062670345
$data1[16,2] = 0.1975827569953026
$data1[16,3] = 0.3745656727807543
$data1[16,4] = 0.2498285737321595
$data1[16,5] = 0.02438596371611723
$data1[17,0] = 0.6442238086993086
$data1[17,1] = 1.748049767009533
$data1[17,2] = 0.1977031169964576
$data1[17,3] = 0.3749232606864226
$data1[17,4] = 0.2503513980993928
$data1[17,5] = 0.02436531028212485
$data1[18,0] = 0.6466974709274567
$data1[18,1] = 1.75212991651361
$data1[18,2] = 0.1971802324417382
$data1[18,3] = 0.3733310579079685
$data1[18,4] = 0.2480168582624094
$data1[18,5] = 0.02445890022088903
$data1[19,0] = 0.6553601761160905
$data1[19,1] = 1.767208947603422
$data1[19,2] = 0.1957527689709053
$data1[19,3] = 0.3682789634530295
$data1[19,4] = 0.2404929073739286
$data1[19,5] = 0.02478462712702267
$data1[20,0] = 0.6612836625660847
$data1[20,1] = 1.778086987525427
$data1[20,2] = 0.1950668228093804
$data1[20,3] = 0.3651992029685758
$data1[20,4] = 0.2358148669995064
$data1[20,5] = 0.02500585565464419
$data1[21,0] = 0.6580982280353709
$data1[21,1] = 1.772187455697917
$data1[21,2] = 0.1954102295218618
$data1[21,3] = 0.3668225865890804
$data1[21,4] = 0.2382896681162698
$data1[21,5] = 0.02488702177371138
$data1[22,0] = 0.6465603064573671
$data1[22,1] = 1.751900703534574
$data1[22,2] = 0.1972077115680477
$data1[22,3] = 0.3734173758145545
$data1[22,4] = 0.248143859443128
$data1[22,5] = 0.02445371813337616
$data1[23,0] = 0.6352398415319058
$data1[23,1] = 1.734398994954972
$data1[23,2] = 0.2001971606101662
$data1[23,3] = 0.3814820952921991
$data1[23,4] = 0.2597943977800323
$data1[23,5] = 0.0241629198062796
$ws.Range("E2:J25").Value = $data1

$data2 = New-Object 'object[,]' 24,1
$data2[0,0] = 0.5908515562150285
$data2[1,0] = 0.5482331789457078
$data2[2,0] = 0.5221342041039918
$data2[3,0] = 0.5115166473494668
$data2[4,0] = 0.5097547162389802
$data2[5,0] = 0.5219909383036878
$data2[6,0] = 0.5761428287554367
$data2[7,0] = 0.6828572811470508
$data2[8,0] = 0.761555235849869
$data2[9,0] = 0.797416490934495
$data2[10,0] = 0.8110044569626069
$data2[11,0] = 0.8080776969133154
$data2[12,0] = 0.7985342230416421
$data2[13,0] = 0.7926896019634171
$data2[14,0] = 0.7592127886031079
$data2[15,0] = 0.7386910077345732
$data2[16,0] = 0.7268932209526753
$data2[17,0] = 0.7228997106452226
$data2[18,0] = 0.7408749913286385
$data2[19,0] = 0.8013371587965423
$data2[20,0] = 0.8408994921579733
$data2[21,0] = 0.819780279826233
$data2[22,0] = 0.7398876108652104
$data2[23,0] = 0.6539347100939779
$ws.Range("M2:M25").Value = $data2

$data3 = New-Object 'object[,]' 24,1
$data3[0,0] = 1.092469087424192
$data3[1,0] = 1.108490097799546
$data3[2,0] = 1.119415655406613
$data3[3,0] = 1.124140962398982
$data3[4,0] = 1.124942067748478
$data3[5,0] = 1.119478277873
$data3[6,0] = 1.097766785633951
$data3[7,0] = 1.063858809494349
$data3[8,0] = 1.04427339248322
$data3[9,0] = 1.036529210048144
$data3[10,0] = 1.033764999486849
$data3[11,0] = 1.034352822341006
$data3[12,0] = 1.036298418529114
$data3[13,0] = 1.037512098453064
$data3[14,0] = 1.044803008313664
$data3[15,0] = 1.049574762500811
$data3[16,0] = 1.052428989078592
$data3[17,0] = 1.053414188321796
$data3[18,0] = 1.049055448213281
$data3[19,0] = 1.035722374853165
$data3[20,0] = 1.02799002889148
$data3[21,0] = 1.03202685831657
$data3[22,0] = 1.049289884871712
$data3[23,0] = 1.072099427495203
$ws.Range("O2:O25").Value = $data3
